$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 1153.9231
$ws.Range("J43").Value = 1230
$ws.Range("L43").Value = 1230
$ws.Range("N43").Value = -1368

# Row 132
$ws.Range("H132").Value = 683.1613
$ws.Range("I132").Value = 692.39343
$ws.Range("J132").Value = 120
$ws.Range("K132").Value = 2077.18029
$ws.Range("L132").Value = 360
$ws.Range("M132").Value = 452.8197100000002
$ws.Range("N132").Value = -5420

# Row 137
$ws.Range("H137").Value = 1740.1364
$ws.Range("I137").Value = 1251.421
$ws.Range("J137").Value = 4835.3335
$ws.Range("K137").Value = 3754.263
$ws.Range("L137").Value = 14506.0005
$ws.Range("M137").Value = -1204.263
$ws.Range("N137").Value = -19606.0005

# Row 138
$ws.Range("H138").Value = 1892.1311
$ws.Range("I138").Value = 1324.5778
$ws.Range("J138").Value = 3488.375
$ws.Range("K138").Value = 3973.7334
$ws.Range("L138").Value = 10465.125
$ws.Range("M138").Value = 1166.2666
$ws.Range("N138").Value = -20745.125

# Row 141
$ws.Range("H141").Value = 7957.023
$ws.Range("I141").Value = 1193.4482
$ws.Range("J141").Value = 21033.268
$ws.Range("K141").Value = 3580.3446
$ws.Range("L141").Value = 63099.804
$ws.Range("M141").Value = 1599.6554
$ws.Range("N141").Value = -73459.804

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4499.8
$ws.Range("I32").Value = 3183.5376
$ws.Range("K32").Value = 3183.5376
$ws.Range("M32").Value = -2896.5376

# Row 132
$ws.Range("H132").Value = 2611.276
$ws.Range("I132").Value = 1512.0526
$ws.Range("K132").Value = 4536.1578
$ws.Range("M132").Value = -2006.1578

$ws = $wb.Worksheets.Item("BSM")
# Row 75
$ws.Range("H75").Value = 13612.5
$ws.Range("I75").Value = 9950
$ws.Range("J75").Value = 14833.333
$ws.Range("K75").Value = 9950
$ws.Range("L75").Value = 14833.333
$ws.Range("M75").Value = -9014
$ws.Range("N75").Value = -16705.333

# Row 78
$ws.Range("H78").Value = 13612.5
$ws.Range("I78").Value = 9950
$ws.Range("J78").Value = 14833.333
$ws.Range("K78").Value = 29850
$ws.Range("L78").Value = 44499.999
$ws.Range("M78").Value = -25170
$ws.Range("N78").Value = -53859.999

# Row 105
$ws.Range("H105").Value = 11367086
$ws.Range("I105").Value = 11367086
$ws.Range("K105").Value = 11367086
$ws.Range("M105").Value = -11365339

# Row 134
$ws.Range("H134").Value = 1989.6842
$ws.Range("I134").Value = 1774
$ws.Range("J134").Value = 2457
$ws.Range("K134").Value = 5322
$ws.Range("L134").Value = 7371
$ws.Range("M134").Value = -2787
$ws.Range("N134").Value = -12441

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2627.0386
$ws.Range("I31").Value = 1720.2413
$ws.Range("J31").Value = 3770.3914
$ws.Range("K31").Value = 1720.2413
$ws.Range("L31").Value = 3770.3914
$ws.Range("M31").Value = -1425.2413
$ws.Range("N31").Value = -4360.3914

# Row 34
$ws.Range("H34").Value = 2627.0386
$ws.Range("I34").Value = 1720.2413
$ws.Range("J34").Value = 3770.3914
$ws.Range("K34").Value = 1720.2413
$ws.Range("L34").Value = 3770.3914
$ws.Range("M34").Value = -1518.2413
$ws.Range("N34").Value = -4174.3914

# Row 50
$ws.Range("H50").Value = 39071.43
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 39071.43
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 39071.43
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -40321.43

# Row 51
$ws.Range("H51").Value = 33099.47
$ws.Range("J51").Value = 33099.47
$ws.Range("L51").Value = 33099.47
$ws.Range("N51").Value = -34571.47

# Row 59
$ws.Range("H59").Value = 39928.57
$ws.Range("J59").Value = 39928.57
$ws.Range("L59").Value = 39928.57
$ws.Range("N59").Value = -42218.57

# Row 60
$ws.Range("H60").Value = 25457.375
$ws.Range("I60").Value = 13133.333
$ws.Range("J60").Value = 28301.385
$ws.Range("K60").Value = 13133.333
$ws.Range("L60").Value = 28301.385
$ws.Range("N60").Value = -29323.385
$ws.Range("M60").Value = -12622.333

# Row 61
$ws.Range("H61").Value = 33099.47
$ws.Range("J61").Value = 33099.47
$ws.Range("L61").Value = 33099.47
$ws.Range("N61").Value = -33795.47

# Row 68
$ws.Range("H68").Value = 24000
$ws.Range("I68").Value = 10000
$ws.Range("J68").Value = 38000
$ws.Range("K68").Value = 10000
$ws.Range("L68").Value = 38000
$ws.Range("N68").Value = -39498
$ws.Range("M68").Value = -9251

# Row 71
$ws.Range("H71").Value = 24000
$ws.Range("I71").Value = 10000
$ws.Range("J71").Value = 38000
$ws.Range("K71").Value = 30000
$ws.Range("L71").Value = 114000
$ws.Range("N71").Value = -121488
$ws.Range("M71").Value = -26256

# Row 132
$ws.Range("H132").Value = 2172.6943
$ws.Range("I132").Value = 704.94446
$ws.Range("J132").Value = 3640.4443
$ws.Range("K132").Value = 2114.83338
$ws.Range("L132").Value = 10921.3329
$ws.Range("M132").Value = 415.16662
$ws.Range("N132").Value = -15981.3329

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 4007.6858
$ws.Range("I131").Value = 548.2857
$ws.Range("J131").Value = 4872.5356
$ws.Range("K131").Value = 1644.8571
$ws.Range("L131").Value = 14617.6068
$ws.Range("M131").Value = 3395.1429
$ws.Range("N131").Value = -24697.6068

$ws = $wb.Worksheets.Item("GSM")
# Row 2
$ws.Range("H2").Value = 1540.7858
$ws.Range("J2").Value = 1256.8889
$ws.Range("L2").Value = 1256.8889
$ws.Range("N2").Value = -1482.8889

# Row 52
$ws.Range("H52").Value = 5509000

# Row 132
$ws.Range("H132").Value = 2481.2593
$ws.Range("I132").Value = 1906.2106
$ws.Range("J132").Value = 3847
$ws.Range("K132").Value = 5718.6318
$ws.Range("L132").Value = 11541
$ws.Range("M132").Value = -3188.6318
$ws.Range("N132").Value = -16601

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 5372.6113
$ws.Range("I132").Value = 5341.5435
$ws.Range("J132").Value = 5427.577
$ws.Range("K132").Value = 16024.6305
$ws.Range("L132").Value = 16282.731
$ws.Range("M132").Value = -13494.6305
$ws.Range("N132").Value = -21342.731

# Row 136
$ws.Range("H136").Value = 14495236
$ws.Range("I136").Value = 3192.9375
$ws.Range("J136").Value = 47619908
$ws.Range("K136").Value = 9578.8125
$ws.Range("L136").Value = 142859724
$ws.Range("M136").Value = -7028.8125
$ws.Range("N136").Value = -142864824

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 49689.332
$ws.Range("I122").Value = 78683.69500000001
$ws.Range("J122").Value = 2573.5
$ws.Range("K122").Value = 236051.085
$ws.Range("L122").Value = 7720.5
$ws.Range("M122").Value = -233601.085
$ws.Range("N122").Value = -12620.5

# Row 132
$ws.Range("H132").Value = 1326.6296
$ws.Range("I132").Value = 1141.6666
$ws.Range("J132").Value = 1807.5333
$ws.Range("K132").Value = 3424.9998
$ws.Range("L132").Value = 5422.5999
$ws.Range("M132").Value = -894.9998000000001
$ws.Range("N132").Value = -10482.5999

# Row 136
$ws.Range("H136").Value = 2866.3333
$ws.Range("I136").Value = 2185.2856
$ws.Range("J136").Value = 3462.25
$ws.Range("K136").Value = 6555.8568
$ws.Range("L136").Value = 10386.75
$ws.Range("M136").Value = -4005.8568
$ws.Range("N136").Value = -15486.75
